$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the Price/Volume columns as text so numeric-looking strings
# (e.g. "1.003", "0.00001084") are not coerced into numbers by Excel.
$valueRange = $ws.Range("D2:E51")
$valueRange.NumberFormat = "@"

$ws.Range("D2").Value = '28.309.87'
$ws.Range("E2").Value = '  +3.16%  '

$ws.Range("D3").Value = '1.817.16'
$ws.Range("E3").Value = '  +4.10%  '

$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '328.40'
$ws.Range("E5").Value = '  +2.03%  '

$ws.Range("E6").Value = '  +0.15%  '

$ws.Range("D7").Value = '0.4352'
$ws.Range("E7").Value = '  +3.56%  '

$ws.Range("D8").Value = '0.3678'
$ws.Range("E8").Value = '  +2.54%  '

$ws.Range("D9").Value = '45.07'
$ws.Range("E9").Value = '  -0.83%  '

$ws.Range("D10").Value = '0.07703'
$ws.Range("E10").Value = '  +3.78%  '

$ws.Range("D11").Value = '1.142'
$ws.Range("E11").Value = '  +2.41%  '

$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.06%  '

$ws.Range("D13").Value = '22.20'
$ws.Range("E13").Value = '  +3.20%  '

$ws.Range("D14").Value = '6.325'
$ws.Range("E14").Value = '  +3.57%  '

$ws.Range("D15").Value = '7.544'
$ws.Range("E15").Value = '  +4.75%  '

$ws.Range("D16").Value = '1.841.15'
$ws.Range("E16").Value = '  +5.98%  '

$ws.Range("D17").Value = '93.19'
$ws.Range("E17").Value = '  +6.37%  '

$ws.Range("D18").Value = '0.00001084'
$ws.Range("E18").Value = '  +1.47%  '

$ws.Range("D19").Value = '0.06530'
$ws.Range("E19").Value = '  +7.59%  '

$ws.Range("E20").Value = '  +0.11%  '

$ws.Range("E21").Value = '  +3.92%  '

$ws.Range("D22").Value = '6.270'
$ws.Range("E22").Value = '  +2.67%  '

$ws.Range("D23").Value = '28.339.41'
$ws.Range("E23").Value = '  +3.17%  '

$ws.Range("E24").Value = '  +1.76%  '

$ws.Range("D25").Value = '2.012'
$ws.Range("E25").Value = '  -14.28%  '

$ws.Range("D26").Value = '162.51'
$ws.Range("E26").Value = '  +7.01%  '

$ws.Range("E27").Value = '  +2.11%  '

$ws.Range("D28").Value = '2.040.36'
$ws.Range("E28").Value = '  +5.32%  '

$ws.Range("D29").Value = '2.296'
$ws.Range("E29").Value = '  -3.74%  '

$ws.Range("D30").Value = '129.14'
$ws.Range("E30").Value = '  +2.72%  '

$ws.Range("D31").Value = '1.223'
$ws.Range("E31").Value = '  +1.92%  '

$ws.Range("D32").Value = '6.008'
$ws.Range("E32").Value = '  +5.49%  '

$ws.Range("D33").Value = '0.09216'
$ws.Range("E33").Value = '  +1.06%  '

$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").Value = '12.97'
$ws.Range("E34").Value = '  +2.02%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '3.394'
$ws.Range("E35").Value = '  -6.50%  '

$ws.Range("D36").Value = '0.02356'
$ws.Range("E36").Value = '  +2.60%  '

$ws.Range("D37").Value = '0.2186'
$ws.Range("E37").Value = '  +2.52%  '

$ws.Range("D38").Value = '5.217'
$ws.Range("E38").Value = '  +2.62%  '

$ws.Range("D39").Value = '0.6591'
$ws.Range("E39").Value = '  +3.03%  '

$ws.Range("D40").Value = '0.06207'
$ws.Range("E40").Value = '  +2.43%  '

$ws.Range("D41").Value = '8.150'
$ws.Range("E41").Value = '  +2.60%  '

$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("D43").Value = '1.443'
$ws.Range("E43").Value = '  +1.35%  '

$ws.Range("E44").Value = '  +0.09%  '

$ws.Range("D45").Value = '14.03'
$ws.Range("E45").Value = '  +2.04%  '

$ws.Range("D46").Value = '0.6124'
$ws.Range("E46").Value = '  +4.60%  '

$ws.Range("D47").Value = '3.754'
$ws.Range("E47").Value = '  +1.14%  '

$ws.Range("D48").Value = '125.94'
$ws.Range("E48").Value = '  +0.18%  '

$ws.Range("D49").Value = '2.024'
$ws.Range("E49").Value = '  +3.73%  '

$ws.Range("D50").Value = '1.160'

$ws.Range("D51").Value = '0.07017'
$ws.Range("E51").Value = '  +2.79%  '

# Restore default (Normal) style on the value range so no stray
# number-format styling is left behind on the cells.
$valueRange.Style = "Normal"
